# Auto-generated edit script: updates crafting-profit figures across all class sheets
# (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) to match refreshed market-board price data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2232.0322
$ws.Range("I15").Value = 2232.0322
$ws.Range("K15").Value = 6696.096600000001
$ws.Range("M15").Value = -6527.096600000001
$ws.Range("H29").Value = 205
$ws.Range("I29").Value = 205
$ws.Range("K29").Value = 615
$ws.Range("M29").Value = -334
$ws.Range("H38").Value = 697.2353000000001
$ws.Range("I38").Value = 57.533333
$ws.Range("J38").Value = 5495
$ws.Range("K38").Value = 172.599999
$ws.Range("L38").Value = 16485
$ws.Range("M38").Value = 199.400001
$ws.Range("N38").Value = -17229
$ws.Range("H43").Value = 1399
$ws.Range("J43").Value = 1148.25
$ws.Range("L43").Value = 1148.25
$ws.Range("N43").Value = -1286.25
$ws.Range("H58").Value = 163
$ws.Range("I58").Value = 163
$ws.Range("K58").Value = 489
$ws.Range("M58").Value = -339
$ws.Range("H64").Value = 10343.728
$ws.Range("J64").Value = 12164.444
$ws.Range("L64").Value = 12164.444
$ws.Range("N64").Value = -12660.444
$ws.Range("H67").Value = 10343.728
$ws.Range("J67").Value = 12164.444
$ws.Range("L67").Value = 12164.444
$ws.Range("N67").Value = -13880.444
$ws.Range("H138").Value = 2848.68
$ws.Range("I138").Value = 911.4706
$ws.Range("J138").Value = 3846.6365
$ws.Range("K138").Value = 2734.4118
$ws.Range("L138").Value = 11539.9095
$ws.Range("M138").Value = 2405.5882
$ws.Range("N138").Value = -21819.9095

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 27780782
$ws.Range("I2").Value = 32260816
$ws.Range("K2").Value = 32260816
$ws.Range("M2").Value = -32260703
$ws.Range("H32").Value = 3512.9387
$ws.Range("I32").Value = 2172.5842
$ws.Range("K32").Value = 2172.5842
$ws.Range("M32").Value = -1885.5842
$ws.Range("H74").Value = 55720.383
$ws.Range("I74").Value = 73509.64
$ws.Range("J74").Value = 20141.857
$ws.Range("K74").Value = 73509.64
$ws.Range("L74").Value = 20141.857
$ws.Range("M74").Value = -72635.64
$ws.Range("N74").Value = -21889.857
$ws.Range("H77").Value = 55720.383
$ws.Range("I77").Value = 73509.64
$ws.Range("J77").Value = 20141.857
$ws.Range("K77").Value = 367548.2
$ws.Range("L77").Value = 100709.285
$ws.Range("M77").Value = -363180.2
$ws.Range("N77").Value = -109445.285
$ws.Range("H116").Value = 27780782
$ws.Range("I116").Value = 32260816
$ws.Range("K116").Value = 32260816
$ws.Range("M116").Value = -32258522
$ws.Range("H122").Value = 2125.8438
$ws.Range("I122").Value = 2076.4285
$ws.Range("J122").Value = 2471.75
$ws.Range("K122").Value = 6229.2855
$ws.Range("L122").Value = 7415.25
$ws.Range("M122").Value = -3779.2855
$ws.Range("N122").Value = -12315.25
$ws.Range("H125").Value = 28888
$ws.Range("J125").Value = 28888
$ws.Range("L125").Value = 28888
$ws.Range("N125").Value = -38728

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 27780782
$ws.Range("I3").Value = 32260816
$ws.Range("K3").Value = 32260816
$ws.Range("M3").Value = -32260702
$ws.Range("H86").Value = 1251.375
$ws.Range("I86").Value = 1251.375
$ws.Range("K86").Value = 1251.375
$ws.Range("M86").Value = -128.375
$ws.Range("H89").Value = 1251.375
$ws.Range("I89").Value = 1251.375
$ws.Range("K89").Value = 6256.875
$ws.Range("M89").Value = -640.875
$ws.Range("H107").Value = 4426.3447
$ws.Range("I107").Value = 1078.0476
$ws.Range("K107").Value = 1078.0476
$ws.Range("M107").Value = 841.9523999999999
$ws.Range("H134").Value = 1982.7709
$ws.Range("I134").Value = 1905.9348
$ws.Range("J134").Value = 3750
$ws.Range("K134").Value = 5717.8044
$ws.Range("L134").Value = 11250
$ws.Range("M134").Value = -3182.8044
$ws.Range("N134").Value = -16320

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2150.0151
$ws.Range("I58").Value = 2007.2632
$ws.Range("K58").Value = 2007.2632
$ws.Range("M58").Value = -1804.2632
$ws.Range("H59").Value = 25089.7
$ws.Range("H136").Value = 2150.0151
$ws.Range("I136").Value = 2007.2632
$ws.Range("K136").Value = 6021.7896
$ws.Range("M136").Value = -3471.7896

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 248.90475
$ws.Range("J12").Value = 292.375
$ws.Range("L12").Value = 877.125
$ws.Range("N12").Value = -1223.125
$ws.Range("H62").Value = 5170.875
$ws.Range("I62").Value = 3341.75
$ws.Range("K62").Value = 10025.25
$ws.Range("M62").Value = -9339.25
$ws.Range("H65").Value = 5170.875
$ws.Range("I65").Value = 3341.75
$ws.Range("K65").Value = 30075.75
$ws.Range("M65").Value = -26643.75
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()
$ws.Range("H107").Value = 1651.4286
$ws.Range("I107").Value = 1030.2
$ws.Range("J107").Value = 1996.5555
$ws.Range("K107").Value = 3090.6
$ws.Range("L107").Value = 5989.666499999999
$ws.Range("M107").Value = -1170.6
$ws.Range("N107").Value = -9829.666499999999
$ws.Range("H131").Value = 37509.215
$ws.Range("I131").Value = 143873.58
$ws.Range("J131").Value = 2054.4285
$ws.Range("K131").Value = 431620.74
$ws.Range("L131").Value = 6163.2855
$ws.Range("M131").Value = -426580.74
$ws.Range("N131").Value = -16243.2855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2035.2142
$ws.Range("I122").Value = 949.375
$ws.Range("J122").Value = 3483
$ws.Range("K122").Value = 2848.125
$ws.Range("L122").Value = 10449
$ws.Range("M122").Value = -398.125
$ws.Range("N122").Value = -15349
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H126").Value = 20838.385
$ws.Range("I126").Value = 25690.1
$ws.Range("K126").Value = 77070.29999999999
$ws.Range("M126").Value = -74600.29999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2157.9473
$ws.Range("I16").Value = 1691.8462
$ws.Range("J16").Value = 3167.8333
$ws.Range("K16").Value = 1691.8462
$ws.Range("L16").Value = 3167.8333
$ws.Range("M16").Value = -1521.8462
$ws.Range("N16").Value = -3507.8333
$ws.Range("H22").Value = 1802.7693
$ws.Range("I22").Value = 554.5
$ws.Range("K22").Value = 554.5
$ws.Range("M22").Value = -259.5
$ws.Range("H27").Value = 1802.7693
$ws.Range("I27").Value = 554.5
$ws.Range("K27").Value = 554.5
$ws.Range("M27").Value = -447.5
$ws.Range("H40").Value = 5221.893
$ws.Range("I40").Value = 4759.6816
$ws.Range("K40").Value = 4759.6816
$ws.Range("M40").Value = -4623.6816
$ws.Range("H61").Value = 33915.766
$ws.Range("I61").Value = 51916.332
$ws.Range("J61").Value = 4837.923
$ws.Range("K61").Value = 51916.332
$ws.Range("L61").Value = 4837.923
$ws.Range("M61").Value = -51714.332
$ws.Range("N61").Value = -5241.923
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H113").Value = 33915.766
$ws.Range("I113").Value = 51916.332
$ws.Range("J113").Value = 4837.923
$ws.Range("K113").Value = 51916.332
$ws.Range("L113").Value = 4837.923
$ws.Range("M113").Value = -49746.332
$ws.Range("N113").Value = -9177.922999999999
$ws.Range("H122").Value = 4771.564
$ws.Range("I122").Value = 4098.375
$ws.Range("J122").Value = 7849
$ws.Range("K122").Value = 12295.125
$ws.Range("L122").Value = 23547
$ws.Range("M122").Value = -9845.125
$ws.Range("N122").Value = -28447
